$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Het Patel"
$ws.Range("B5").Value = "hetpatel5542@gmail.com"
$ws.Range("C5").Value = "GCET"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "7698545581"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "STATIC_COMBO"
$ws.Range("F5").Value = "OFFLINE"
